$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the "_old"/"_new" header suffixes to "_FV2210"/"_FV2304" -------
# This touches the 20 header cells in row 1 (A1:J1 = "_old" columns,
# L1:U1 = "_new" columns; K1 = "diff" is untouched).
$ws.Cells.Replace("_old", "_FV2210")
$ws.Cells.Replace("_new", "_FV2304")

# --- Turn the data range into an Excel Table (ListObject) -----------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U66"), $null, 1)
$tbl.Name = "Table1"

# --- Freeze the header row --------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
